$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.1953125
$ws.Cells.Item(2, 3).Value = 0.52734375
$ws.Cells.Item(2, 10).Value = 0.0390625
$ws.Cells.Item(2, 16).Value = 0.1171875
$ws.Cells.Item(2, 19).Value = 0.12109375

# Row 3
$ws.Cells.Item(3, 2).Value = 0.01360544217687075
$ws.Cells.Item(3, 3).Value = 0.034013605442176867
$ws.Cells.Item(3, 10).Value = 0.054421768707482991
$ws.Cells.Item(3, 16).Value = 0.73469387755102045
$ws.Cells.Item(3, 19).Value = 0.16326530612244899

# Row 4
$ws.Cells.Item(4, 10).Value = 0.088235294117647065
$ws.Cells.Item(4, 16).Value = 0.70588235294117652
$ws.Cells.Item(4, 19).Value = 0.20588235294117649

# Row 6
$ws.Cells.Item(6, 2).Value = 0.072072072072072071
$ws.Cells.Item(6, 4).Value = 0.0045045045045045036
$ws.Cells.Item(6, 6).Value = 0.036036036036036043
$ws.Cells.Item(6, 10).Value = 0.26126126126126131
$ws.Cells.Item(6, 15).Value = 0.027027027027027029
$ws.Cells.Item(6, 17).Value = 0.14864864864864871
$ws.Cells.Item(6, 18).Value = 0.04954954954954955
$ws.Cells.Item(6, 19).Value = 0.40090090090090091

# Row 7
$ws.Cells.Item(7, 2).Value = 0.064285714285714279
$ws.Cells.Item(7, 4).Value = 0.01428571428571429
$ws.Cells.Item(7, 6).Value = 0.050000000000000003
$ws.Cells.Item(7, 10).Value = 0.16428571428571431
$ws.Cells.Item(7, 17).Value = 0.1714285714285714
$ws.Cells.Item(7, 18).Value = 0.071428571428571425
$ws.Cells.Item(7, 19).Value = 0.4642857142857143

# Row 8
$ws.Cells.Item(8, 2).Value = 0.083532219570405727
$ws.Cells.Item(8, 4).Value = 0.021479713603818611
$ws.Cells.Item(8, 6).Value = 0.066825775656324582
$ws.Cells.Item(8, 10).Value = 0.095465393794749401
$ws.Cells.Item(8, 15).Value = 0.02386634844868735
$ws.Cells.Item(8, 17).Value = 0.17183770883054891
$ws.Cells.Item(8, 18).Value = 0.083532219570405727
$ws.Cells.Item(8, 19).Value = 0.45346062052505959

# Row 9
$ws.Cells.Item(9, 2).Value = 0.054545454545454543
$ws.Cells.Item(9, 4).Value = 0.018181818181818181
$ws.Cells.Item(9, 5).Value = 0.0036363636363636359
$ws.Cells.Item(9, 6).Value = 0.080000000000000002
$ws.Cells.Item(9, 10).Value = 0.10181818181818179
$ws.Cells.Item(9, 15).Value = 0.014545454545454551
$ws.Cells.Item(9, 17).Value = 0.21090909090909091
$ws.Cells.Item(9, 18).Value = 0.08727272727272728
$ws.Cells.Item(9, 19).Value = 0.42909090909090908

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1217228464419476
$ws.Cells.Item(10, 4).Value = 0.017790262172284639
$ws.Cells.Item(10, 6).Value = 0.084269662921348312
$ws.Cells.Item(10, 10).Value = 0.12453183520599249
$ws.Cells.Item(10, 15).Value = 0.016853932584269659
$ws.Cells.Item(10, 17).Value = 0.2134831460674157
$ws.Cells.Item(10, 18).Value = 0.066479400749063666
$ws.Cells.Item(10, 19).Value = 0.35486891385767788

# Row 11
$ws.Cells.Item(11, 7).Value = 0.1133004926108374
$ws.Cells.Item(11, 10).Value = 0.13300492610837439
$ws.Cells.Item(11, 11).Value = 0.14778325123152711
$ws.Cells.Item(11, 12).Value = 0.59605911330049266
$ws.Cells.Item(11, 19).Value = 0.009852216748768473

# Row 12
$ws.Cells.Item(12, 7).Value = 0.74045801526717558
$ws.Cells.Item(12, 10).Value = 0.12977099236641221
$ws.Cells.Item(12, 11).Value = 0.045801526717557252
$ws.Cells.Item(12, 12).Value = 0.068702290076335881
$ws.Cells.Item(12, 19).Value = 0.01526717557251908

# Row 13
$ws.Cells.Item(13, 7).Value = 0.69999999999999996
$ws.Cells.Item(13, 10).Value = 0.23333333333333331
$ws.Cells.Item(13, 19).Value = 0.066666666666666666

# Row 15
$ws.Cells.Item(15, 6).Value = 0.030927835051546389
$ws.Cells.Item(15, 8).Value = 0.16494845360824739
$ws.Cells.Item(15, 9).Value = 0.10309278350515461
$ws.Cells.Item(15, 10).Value = 0.35051546391752569
$ws.Cells.Item(15, 11).Value = 0.041237113402061848
$ws.Cells.Item(15, 14).Value = 0.0051546391752577319
$ws.Cells.Item(15, 15).Value = 0.03608247422680412
$ws.Cells.Item(15, 19).Value = 0.26804123711340211

# Row 16
$ws.Cells.Item(16, 6).Value = 0.050632911392405063
$ws.Cells.Item(16, 8).Value = 0.120253164556962
$ws.Cells.Item(16, 9).Value = 0.10759493670886081
$ws.Cells.Item(16, 10).Value = 0.34177215189873422
$ws.Cells.Item(16, 11).Value = 0.1012658227848101
$ws.Cells.Item(16, 13).Value = 0.0063291139240506328
$ws.Cells.Item(16, 14).Value = 0.012658227848101271
$ws.Cells.Item(16, 15).Value = 0.031645569620253167
$ws.Cells.Item(16, 19).Value = 0.22784810126582281

# Row 17
$ws.Cells.Item(17, 6).Value = 0.01210653753026634
$ws.Cells.Item(17, 8).Value = 0.1912832929782082
$ws.Cells.Item(17, 9).Value = 0.12832929782082331
$ws.Cells.Item(17, 10).Value = 0.38740920096852299
$ws.Cells.Item(17, 11).Value = 0.082324455205811137
$ws.Cells.Item(17, 13).Value = 0.01210653753026634
$ws.Cells.Item(17, 14).Value = 0.0024213075060532689
$ws.Cells.Item(17, 15).Value = 0.050847457627118647
$ws.Cells.Item(17, 19).Value = 0.13317191283292981

# Row 18
$ws.Cells.Item(18, 6).Value = 0.032894736842105261
$ws.Cells.Item(18, 8).Value = 0.19736842105263161
$ws.Cells.Item(18, 9).Value = 0.1118421052631579
$ws.Cells.Item(18, 10).Value = 0.34868421052631582
$ws.Cells.Item(18, 11).Value = 0.1118421052631579
$ws.Cells.Item(18, 13).Value = 0.0065789473684210523
$ws.Cells.Item(18, 15).Value = 0.078947368421052627
$ws.Cells.Item(18, 19).Value = 0.1118421052631579

# Row 19
$ws.Cells.Item(19, 6).Value = 0.01377633711507293
$ws.Cells.Item(19, 8).Value = 0.2074554294975689
$ws.Cells.Item(19, 9).Value = 0.13533225283630471
$ws.Cells.Item(19, 10).Value = 0.3354943273905997
$ws.Cells.Item(19, 11).Value = 0.075364667747163702
$ws.Cells.Item(19, 13).Value = 0.0186385737439222
$ws.Cells.Item(19, 14).Value = 0.002431118314424636
$ws.Cells.Item(19, 15).Value = 0.072123176661264179
$ws.Cells.Item(19, 19).Value = 0.13938411669367909
